# Update cryptocurrency price/volume data per latest GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.524.98'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.05%  '
$ws.Range('E2').NumberFormat = 'General'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.605.09'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.06%  '
$ws.Range('E3').NumberFormat = 'General'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E4').NumberFormat = 'General'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.16'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('E5').NumberFormat = 'General'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.519'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +6.96%  '
$ws.Range('E6').NumberFormat = 'General'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E7').NumberFormat = 'General'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '27.01'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +9.62%  '
$ws.Range('E8').NumberFormat = 'General'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.42'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('E9').NumberFormat = 'General'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.97%  '
$ws.Range('E10').NumberFormat = 'General'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0599'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.55%  '
$ws.Range('E11').NumberFormat = 'General'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0911'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('E12').NumberFormat = 'General'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.834.01'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.99%  '
$ws.Range('E13').NumberFormat = 'General'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.622.43'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.32%  '
$ws.Range('E14').NumberFormat = 'General'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.559.93'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.22%  '
$ws.Range('E15').NumberFormat = 'General'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.537'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.60%  '
$ws.Range('E16').NumberFormat = 'General'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.61%  '
$ws.Range('E17').NumberFormat = 'General'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.57'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.62%  '
$ws.Range('E18').NumberFormat = 'General'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.61'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.74%  '
$ws.Range('E19').NumberFormat = 'General'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.63'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.21%  '
$ws.Range('E20').NumberFormat = 'General'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.92%  '
$ws.Range('E21').NumberFormat = 'General'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('E22').NumberFormat = 'General'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.71%  '
$ws.Range('E23').NumberFormat = 'General'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.17'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.06%  '
$ws.Range('E24').NumberFormat = 'General'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.08'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.29%  '
$ws.Range('E25').NumberFormat = 'General'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.34'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.65%  '
$ws.Range('E26').NumberFormat = 'General'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.32'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.67%  '
$ws.Range('E27').NumberFormat = 'General'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.109'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.60%  '
$ws.Range('E28').NumberFormat = 'General'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.40'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.59%  '
$ws.Range('E29').NumberFormat = 'General'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E30').NumberFormat = 'General'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0472'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.81%  '
$ws.Range('E31').NumberFormat = 'General'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('E32').NumberFormat = 'General'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.66%  '
$ws.Range('E33').NumberFormat = 'General'
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('B34').NumberFormat = 'General'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C34').NumberFormat = 'General'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.10'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.56%  '
$ws.Range('E34').NumberFormat = 'General'
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Maker'
$ws.Range('B35').NumberFormat = 'General'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C35').NumberFormat = 'General'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.418.06'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.66%  '
$ws.Range('E35').NumberFormat = 'General'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.02'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('E36').NumberFormat = 'General'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.52'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.18%  '
$ws.Range('E37').NumberFormat = 'General'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.79'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +5.30%  '
$ws.Range('E38').NumberFormat = 'General'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.20%  '
$ws.Range('E39').NumberFormat = 'General'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.20%  '
$ws.Range('E40').NumberFormat = 'General'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.534'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.18%  '
$ws.Range('E41').NumberFormat = 'General'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.16%  '
$ws.Range('E42').NumberFormat = 'General'
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'BitcoinSV'
$ws.Range('B43').NumberFormat = 'General'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('C43').NumberFormat = 'General'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '52.95'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +21.17%  '
$ws.Range('E43').NumberFormat = 'General'
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('B44').NumberFormat = 'General'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('C44').NumberFormat = 'General'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E44').NumberFormat = 'General'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.793'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.25%  '
$ws.Range('E45').NumberFormat = 'General'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.00%  '
$ws.Range('E46').NumberFormat = 'General'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.73'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.49%  '
$ws.Range('E47').NumberFormat = 'General'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.29'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E48').NumberFormat = 'General'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.746.95'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.21%  '
$ws.Range('E49').NumberFormat = 'General'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('E50').NumberFormat = 'General'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.08%  '
$ws.Range('E51').NumberFormat = 'General'
